$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the new column header values in row 15 (AH15:BB15)
$ws.Range("AH15").Value = "biotic_relationship"
$ws.Range("AI15").Value = "dew_point"
$ws.Range("AJ15").Value = "extrachrom_elements"
$ws.Range("AK15").Value = "indoor_surf"
$ws.Range("AL15").Value = "isolation_source"
$ws.Range("AM15").Value = "rel_to_oxygen"
$ws.Range("AN15").Value = "samp_collect_device"
$ws.Range("AO15").Value = "samp_mat_process"
$ws.Range("AP15").Value = "samp_size"
$ws.Range("AQ15").Value = "samp_sort_meth"
$ws.Range("AR15").Value = "samp_vol_we_dna_ext"
$ws.Range("AS15").Value = "source_material_id"
$ws.Range("AT15").Value = "subspecf_gen_lin"
$ws.Range("AU15").Value = "substructure_type"
$ws.Range("AV15").Value = "surf_air_cont"
$ws.Range("AW15").Value = "surf_humidity"
$ws.Range("AX15").Value = "surf_material"
$ws.Range("AY15").Value = "surf_moisture"
$ws.Range("AZ15").Value = "surf_moisture_ph"
$ws.Range("BA15").Value = "surf_temp"
$ws.Range("BB15").Value = "trophic_level"

# Copy the "optional" (yellow) header style from C15 onto the new header cells
$ws.Range("C15").Copy() | Out-Null
$ws.Range("AH15:BB15").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Add the field-definition cell comments
$ws.Range("AH15").AddComment("Free-living or from host (define relationship)") | Out-Null
$ws.Range("AI15").AddComment("temperature to which a given parcel of humid air must be cooled, at constant barometric pressure, for water vapor to condense into water.") | Out-Null
$ws.Range("AJ15").AddComment("Plasmids that have significance phenotypic consequence") | Out-Null
$ws.Range("AK15").AddComment("type of indoor surface") | Out-Null
$ws.Range("AL15").AddComment("Describes the physical, environmental and/or local geographical source of the biological sample from which the sample was derived.") | Out-Null
$ws.Range("AM15").AddComment("Aerobic or anaerobic") | Out-Null
$ws.Range("AN15").AddComment("Method or device employed for collecting sample") | Out-Null
$ws.Range("AO15").AddComment("Processing applied to the sample during or after isolation") | Out-Null
$ws.Range("AP15").AddComment("Amount or size of sample (volume, mass or area) that was collected") | Out-Null
$ws.Range("AQ15").AddComment("method by which samples are sorted") | Out-Null
$ws.Range("AR15").AddComment("volume (mL) or weight (g) of sample processed for DNA extraction") | Out-Null
$ws.Range("AS15").AddComment("unique identifier assigned to a material sample used for extracting nucleic acids, and subsequent sequencing. The identifier can refer either to the original material collected or to any derived sub-samples.") | Out-Null
$ws.Range("AT15").AddComment("Information about the genetic distinctness of the lineage (eg., biovar, serovar)") | Out-Null
$ws.Range("AU15").AddComment("substructure or under building is that largely hidden section of the building which is built off the foundations to the ground floor level") | Out-Null
$ws.Range("AV15").AddComment("contaminant identified on surface") | Out-Null
$ws.Range("AW15").AddComment("surfaces: water activity as a function of air and material moisture") | Out-Null
$ws.Range("AX15").AddComment("surface materials at the point of sampling") | Out-Null
$ws.Range("AY15").AddComment("water held on a surface") | Out-Null
$ws.Range("AZ15").AddComment("pH measurement of surface") | Out-Null
$ws.Range("BA15").AddComment("temperature of the surface at the time of sampling") | Out-Null
$ws.Range("BB15").AddComment("Feeding position in food chain (eg., chemolithotroph)") | Out-Null
